$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "335.27"
Set-TextValue $ws "E2" "1.82%"
Set-TextValue $ws "G2" "10"
Set-TextValue $ws "D3" "43.80"
Set-TextValue $ws "E3" "6.29%"
Set-TextValue $ws "G3" "10"
Set-TextValue $ws "D4" "5.758"
Set-TextValue $ws "E4" "2.20%"
Set-TextValue $ws "G4" "10"
Set-TextValue $ws "D5" "0.08322"
Set-TextValue $ws "E5" "1.61%"
Set-TextValue $ws "G5" "10"
Set-TextValue $ws "D6" "8.841"
Set-TextValue $ws "E6" "1.09%"
Set-TextValue $ws "G6" "10"
Set-TextValue $ws "B7" "FTXToken"
Set-TextValue $ws "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D7" "1.965"
Set-TextValue $ws "E7" "-2.35%"
Set-TextValue $ws "G7" "10"
Set-TextValue $ws "B8" "BTSEToken"
Set-TextValue $ws "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D8" "2.889"
Set-TextValue $ws "E8" "-2.56%"
Set-TextValue $ws "G8" "10"
Set-TextValue $ws "B9" "MXToken"
Set-TextValue $ws "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D9" "0.9427"
Set-TextValue $ws "E9" "2.49%"
Set-TextValue $ws "G9" "10"
Set-TextValue $ws "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D10" "0.1245"
Set-TextValue $ws "E10" "-2.93%"
Set-TextValue $ws "G10" "10"
Set-TextValue $ws "B11" "WazirX"
Set-TextValue $ws "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D11" "0.1969"
Set-TextValue $ws "E11" "0.92%"
Set-TextValue $ws "G11" "10"
Set-TextValue $ws "B12" "MandalaExchangeToken"
Set-TextValue $ws "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D12" "0.1024"
Set-TextValue $ws "E12" "9.22%"
Set-TextValue $ws "G12" "10"
Set-TextValue $ws "B13" "BitrueCoin"
Set-TextValue $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D13" "0.04573"
Set-TextValue $ws "E13" "17.50%"
Set-TextValue $ws "G13" "10"
Set-TextValue $ws "B14" "BitMartToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D14" "0.1068"
Set-TextValue $ws "E14" "0.75%"
Set-TextValue $ws "G14" "10"
Set-TextValue $ws "B15" "BitForexToken"
Set-TextValue $ws "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D15" "0.001298"
Set-TextValue $ws "E15" "-0.45%"
Set-TextValue $ws "G15" "10"
Set-TextValue $ws "B16" "TigerCash"
Set-TextValue $ws "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D16" "0.005964"
Set-TextValue $ws "E16" "-5.14%"
Set-TextValue $ws "G16" "10"
Set-TextValue $ws "B17" "LEO"
Set-TextValue $ws "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.501"
Set-TextValue $ws "E17" "1.54%"
Set-TextValue $ws "G17" "10"
Set-TextValue $ws "B18" "GateToken"
Set-TextValue $ws "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D18" "4.522"
Set-TextValue $ws "E18" "0.65%"
Set-TextValue $ws "G18" "10"
Set-TextValue $ws "G19" "10"
Set-TextValue $ws "D20" "8.745"
Set-TextValue $ws "E20" "6.21%"
Set-TextValue $ws "G20" "10"
Set-TextValue $ws "D21" "0.1364"
Set-TextValue $ws "E21" "-0.13%"
Set-TextValue $ws "G21" "10"
Set-TextValue $ws "E22" "11.80%"
Set-TextValue $ws "G22" "10"
Set-TextValue $ws "D23" "0.04418"
Set-TextValue $ws "E23" "0.54%"
Set-TextValue $ws "G23" "10"
Set-TextValue $ws "D24" "0.001263"
Set-TextValue $ws "E24" "0.60%"
Set-TextValue $ws "G24" "10"
Set-TextValue $ws "D25" "0.004349"
Set-TextValue $ws "E25" "0.85%"
Set-TextValue $ws "G25" "10"
Set-TextValue $ws "D26" "0.0001262"
Set-TextValue $ws "E26" "5.14%"
Set-TextValue $ws "G26" "10"
Set-TextValue $ws "D27" "0.0003997"
Set-TextValue $ws "G27" "10"
Set-TextValue $ws "G28" "10"
Set-TextValue $ws "G29" "10"
Set-TextValue $ws "G30" "10"
Set-TextValue $ws "G31" "10"
Set-TextValue $ws "G32" "10"
Set-TextValue $ws "G33" "10"
Set-TextValue $ws "G34" "10"
Set-TextValue $ws "G35" "10"
Set-TextValue $ws "G36" "10"
Set-TextValue $ws "G37" "10"
Set-TextValue $ws "G38" "10"
Set-TextValue $ws "D39" "0.02804"
Set-TextValue $ws "E39" "0.87%"
Set-TextValue $ws "G39" "10"
Set-TextValue $ws "D40" "0.06097"
Set-TextValue $ws "E40" "12.83%"
Set-TextValue $ws "G40" "10"
Set-TextValue $ws "D41" "0.007915"
Set-TextValue $ws "E41" "1.52%"
Set-TextValue $ws "G41" "10"
Set-TextValue $ws "E42" "0.81%"
Set-TextValue $ws "G42" "10"
Set-TextValue $ws "D43" "0.008982"
Set-TextValue $ws "E43" "0.42%"
Set-TextValue $ws "G43" "10"
Set-TextValue $ws "D44" "0.002144"
Set-TextValue $ws "E44" "-1.25%"
Set-TextValue $ws "G44" "10"
Set-TextValue $ws "D45" "0.01039"
Set-TextValue $ws "E45" "-13.74%"
Set-TextValue $ws "G45" "10"
Set-TextValue $ws "D46" "0.00007014"
Set-TextValue $ws "E46" "3.52%"
Set-TextValue $ws "G46" "10"
Set-TextValue $ws "D47" "0.00000000751"
Set-TextValue $ws "E47" "0.15%"
Set-TextValue $ws "G47" "10"
Set-TextValue $ws "E48" "0.01%"
Set-TextValue $ws "G48" "10"
Set-TextValue $ws "E49" "-0.25%"
Set-TextValue $ws "G49" "10"
Set-TextValue $ws "E50" "0.15%"
Set-TextValue $ws "G50" "10"
Set-TextValue $ws "E51" "0.15%"
Set-TextValue $ws "G51" "10"
